$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 1160996
$ws.Range("C4").Value = 222
$ws.Range("D4").Value = 173725
$ws.Range("E4").Value = 919823
$ws.Range("F4").Value = 16475
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 67448

$ws.Range("A30").Value = "Bielorrusia"
$ws.Range("B30").Value = 16705
$ws.Range("C30").Value = 877
$ws.Range("D30").Value = 3196
$ws.Range("E30").Value = 13410
$ws.Range("F30").Value = 92
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 99

$ws.Range("A31").Value = "Israel"
$ws.Range("B31").Value = 16193
$ws.Range("C31").Value = 8
$ws.Range("D31").Value = 9634
$ws.Range("E31").Value = 6329
$ws.Range("F31").Value = 103
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 230

$ws.Range("A50").Value = "Australia"
$ws.Range("B50").Value = 6801
$ws.Range("C50").Value = 20
$ws.Range("D50").Value = 5817
$ws.Range("E50").Value = 889
$ws.Range("F50").Value = 29
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 95

$ws.Range("A54").Value = "Finlandia"
$ws.Range("B54").Value = 5254
$ws.Range("C54").Value = 78
$ws.Range("D54").Value = 3000
$ws.Range("E54").Value = 2024
$ws.Range("F54").Value = 49
$ws.Range("G54").Value = 10
$ws.Range("H54").Value = 230

$ws.Range("A62").Value = "Barein"
$ws.Range("B62").Value = 3356
$ws.Range("C62").Value = 72
$ws.Range("D62").Value = 1717
$ws.Range("E62").Value = 1631
$ws.Range("F62").Value = 1
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 8

$ws.Range("A99").Value = "Libano"
$ws.Range("B99").Value = 737
$ws.Range("C99").Value = 4
$ws.Range("D99").Value = 200
$ws.Range("E99").Value = 512
$ws.Range("F99").Value = 43
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 25

$ws.Range("A100").Value = "Niger"
$ws.Range("B100").Value = 736
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 507
$ws.Range("E100").Value = 194
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 35

$ws.Range("A111").Value = "San Marino"
$ws.Range("B111").Value = 582
$ws.Range("C111").Value = 2
$ws.Range("D111").Value = 86
$ws.Range("E111").Value = 455
$ws.Range("F111").Value = 5
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 41

$ws.Range("A117").Value = "Malta"
$ws.Range("B117").Value = 477
$ws.Range("C117").Value = 9
$ws.Range("D117").Value = 392
$ws.Range("E117").Value = 81
$ws.Range("F117").Value = 1
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 4

$ws.Range("A142").Value = "Madagascar"
$ws.Range("B142").Value = 149
$ws.Range("C142").Value = 14
$ws.Range("D142").Value = 98
$ws.Range("E142").Value = 51
$ws.Range("F142").Value = 1
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 0

$ws.Range("A143").Value = "Gibraltar"
$ws.Range("B143").Value = 144
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 132
$ws.Range("E143").Value = 12
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 0

$ws.Range("A144").Value = "Brunei"
$ws.Range("B144").Value = 138
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 126
$ws.Range("E144").Value = 11
$ws.Range("F144").Value = 2
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 1

$ws.Range("A145").Value = "Etiopia"
$ws.Range("B145").Value = 135
$ws.Range("C145").Value = 2
$ws.Range("D145").Value = 75
$ws.Range("E145").Value = 57
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 3

$ws.Range("A168").Value = "Nepal"
$ws.Range("B168").Value = 69
$ws.Range("C168").Value = 10
$ws.Range("D168").Value = 16
$ws.Range("E168").Value = 53
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 0

$ws.Range("A169").Value = "Libia"
$ws.Range("B169").Value = 63
$ws.Range("C169").Value = 0
$ws.Range("D169").Value = 22
$ws.Range("E169").Value = 38
$ws.Range("F169").Value = 0
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 3

$ws.Range("A172").Value = "Macao"
$ws.Range("B172").Value = 45
$ws.Range("C172").Value = 0
$ws.Range("D172").Value = 39
$ws.Range("E172").Value = 6
$ws.Range("F172").Value = 1
$ws.Range("G172").Value = 0
$ws.Range("H172").Value = 0

$ws.Range("A198").Value = "San Cristobal y Nieves"
$ws.Range("B198").Value = 15
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 8
$ws.Range("E198").Value = 7
$ws.Range("F198").Value = 0
$ws.Range("G198").Value = 0
$ws.Range("H198").Value = 0

$ws.Range("A199").Value = "Burundi"
$ws.Range("B199").Value = 15
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 7
$ws.Range("E199").Value = 7
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 1

